# Update "想去人数" (F column) values across the workbook sheets to reflect
# newly generated output (commit: Update gh-pages to output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1359
$ws1.Range("F5").Value  = 250
$ws1.Range("F7").Value  = 986
$ws1.Range("F8").Value  = 19076
$ws1.Range("F13").Value = 375
$ws1.Range("F22").Value = 143

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value  = 115
$ws2.Range("F16").Value = 78

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 608
$ws3.Range("F4").Value = 571

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 608
$ws4.Range("F5").Value  = 571
$ws4.Range("F8").Value  = 1359
$ws4.Range("F10").Value = 250
$ws4.Range("F14").Value = 986
$ws4.Range("F15").Value = 19076
$ws4.Range("F19").Value = 115
$ws4.Range("F26").Value = 375
$ws4.Range("F37").Value = 78
$ws4.Range("F39").Value = 143
